$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")

# ALC!row15
$ws1.Range("H15").Value = 1914.4487
$ws1.Range("I15").Value = 1914.4487
$ws1.Range("K15").Value = 5743.3461
$ws1.Range("M15").Value = -5574.3461

# ALC!row96
$ws1.Range("H96").Value = 27778390
$ws1.Range("I96").Value = 41667196
$ws1.Range("K96").Value = 125001588
$ws1.Range("M96").Value = -125000215

# ALC!row116
$ws1.Range("H116").Value = 4794.8237
$ws1.Range("I116").Value = 2357.1428
$ws1.Range("K116").Value = 2357.1428
$ws1.Range("M116").Value = 1084.8572

# ALC!row129
$ws1.Range("H129").Value = 1110.25
$ws1.Range("J129").Value = 1228.1842
$ws1.Range("L129").Value = 3684.5526
$ws1.Range("N129").Value = -13684.5526

# ALC!row132
$ws1.Range("H132").Value = 3169.658
$ws1.Range("I132").Value = 3091.6667
$ws1.Range("J132").Value = 3462.125
$ws1.Range("K132").Value = 9275.000100000001
$ws1.Range("L132").Value = 10386.375
$ws1.Range("M132").Value = -6745.000100000001
$ws1.Range("N132").Value = -15446.375

# ARM!row32
$ws2.Range("H32").Value = 4845.983
$ws2.Range("I32").Value = 3740.1
$ws2.Range("J32").Value = 11757.75
$ws2.Range("K32").Value = 3740.1
$ws2.Range("L32").Value = 11757.75
$ws2.Range("M32").Value = -3453.1
$ws2.Range("N32").Value = -12331.75

# ARM!row45
$ws2.Range("H45").Value = 4055.818
$ws2.Range("I45").Value = 3766.6667
$ws2.Range("J45").Value = 4402.8
$ws2.Range("K45").Value = 3766.6667
$ws2.Range("L45").Value = 4402.8
$ws2.Range("M45").Value = -3389.6667
$ws2.Range("N45").Value = -5156.8

# ARM!row132
$ws2.Range("H132").Value = 18542.613
$ws2.Range("I132").Value = 2229.0435
$ws2.Range("J132").Value = 65444.125
$ws2.Range("K132").Value = 6687.130500000001
$ws2.Range("L132").Value = 196332.375
$ws2.Range("M132").Value = -4157.130500000001
$ws2.Range("N132").Value = -201392.375

# BSM!row107
$ws3.Range("H107").Value = 784.6667
$ws3.Range("J107").Value = 746.3333
$ws3.Range("L107").Value = 746.3333
$ws3.Range("N107").Value = -4586.3333

# BSM!row113
$ws3.Range("H113").Value = 0
$ws3.Range("I113").Value = 0
$ws3.Range("K113").Value = 0
$ws3.Range("M113").ClearContents()

# CRP!row31
$ws4.Range("H31").Value = 3603.4856
$ws4.Range("I31").Value = 3676.2222
$ws4.Range("J31").Value = 3578.3076
$ws4.Range("K31").Value = 3676.2222
$ws4.Range("L31").Value = 3578.3076
$ws4.Range("M31").Value = -3381.2222
$ws4.Range("N31").Value = -4168.3076

# CRP!row33
$ws4.Range("H33").Value = 2031
$ws4.Range("I33").Value = 2031
$ws4.Range("K33").Value = 2031
$ws4.Range("M33").Value = -1652

# CRP!row34
$ws4.Range("H34").Value = 3603.4856
$ws4.Range("I34").Value = 3676.2222
$ws4.Range("J34").Value = 3578.3076
$ws4.Range("K34").Value = 3676.2222
$ws4.Range("L34").Value = 3578.3076
$ws4.Range("M34").Value = -3474.2222
$ws4.Range("N34").Value = -3982.3076

# CRP!row39
$ws4.Range("H39").Value = 0
$ws4.Range("I39").Value = 0
$ws4.Range("K39").Value = 0
$ws4.Range("M39").ClearContents()

# CRP!row41
$ws4.Range("H41").Value = 0
$ws4.Range("I41").Value = 0
$ws4.Range("K41").Value = 0
$ws4.Range("M41").ClearContents()

# CRP!row49
$ws4.Range("H49").Value = 0
$ws4.Range("I49").Value = 0
$ws4.Range("K49").Value = 0
$ws4.Range("M49").ClearContents()

# CRP!row99
$ws4.Range("H99").Value = 41671230
$ws4.Range("I99").Value = 3428.5715
$ws4.Range("J99").Value = 100006160
$ws4.Range("K99").Value = 3428.5715
$ws4.Range("L99").Value = 100006160
$ws4.Range("M99").Value = -1930.5715
$ws4.Range("N99").Value = -100009156

# CRP!row126
$ws4.Range("H126").Value = 41671230
$ws4.Range("I126").Value = 3428.5715
$ws4.Range("J126").Value = 100006160
$ws4.Range("K126").Value = 10285.7145
$ws4.Range("L126").Value = 300018480
$ws4.Range("M126").Value = -7815.7145
$ws4.Range("N126").Value = -300023420

# CRP!row132
$ws4.Range("H132").Value = 3723.3333
$ws4.Range("I132").Value = 2934.4614
$ws4.Range("J132").Value = 5005.25
$ws4.Range("K132").Value = 8803.3842
$ws4.Range("L132").Value = 15015.75
$ws4.Range("M132").Value = -6273.3842
$ws4.Range("N132").Value = -20075.75

# CUL!row131
$ws5.Range("H131").Value = 701.45
$ws5.Range("I131").Value = 431.42856
$ws5.Range("J131").Value = 721.7742
$ws5.Range("K131").Value = 1294.28568
$ws5.Range("L131").Value = 2165.3226
$ws5.Range("M131").Value = 3745.71432
$ws5.Range("N131").Value = -12245.3226

# GSM!row80
$ws6.Range("H80").Value = 3682.3704
$ws6.Range("I80").Value = 2744.3333
$ws6.Range("J80").Value = 4151.3887
$ws6.Range("K80").Value = 2744.3333
$ws6.Range("L80").Value = 4151.3887
$ws6.Range("M80").Value = -1746.3333
$ws6.Range("N80").Value = -6147.3887

# GSM!row83
$ws6.Range("H83").Value = 3682.3704
$ws6.Range("I83").Value = 2744.3333
$ws6.Range("J83").Value = 4151.3887
$ws6.Range("K83").Value = 13721.6665
$ws6.Range("L83").Value = 20756.9435
$ws6.Range("M83").Value = -8729.666499999999
$ws6.Range("N83").Value = -30740.9435

# GSM!row102
$ws6.Range("H102").Value = 27780930
$ws6.Range("I102").Value = 35717396
$ws6.Range("J102").Value = 3298.5
$ws6.Range("K102").Value = 35717396
$ws6.Range("L102").Value = 3298.5
$ws6.Range("M102").Value = -35715774
$ws6.Range("N102").Value = -6542.5

# GSM!row113
$ws6.Range("H113").Value = 2213.45
$ws6.Range("I113").Value = 1635.0769
$ws6.Range("J113").Value = 3287.5715
$ws6.Range("K113").Value = 1635.0769
$ws6.Range("L113").Value = 3287.5715
$ws6.Range("M113").Value = 534.9231
$ws6.Range("N113").Value = -7627.5715

# GSM!row122
$ws6.Range("H122").Value = 111112940
$ws6.Range("I122").Value = 37039090
$ws6.Range("J122").Value = 333334500
$ws6.Range("K122").Value = 111117270
$ws6.Range("L122").Value = 1000003500
$ws6.Range("M122").Value = -111114820
$ws6.Range("N122").Value = -1000008400

# GSM!row126
$ws6.Range("H126").Value = 5536.5454
$ws6.Range("I126").Value = 4500
$ws6.Range("J126").Value = 6780.4
$ws6.Range("K126").Value = 13500
$ws6.Range("L126").Value = 20341.2
$ws6.Range("M126").Value = -11030
$ws6.Range("N126").Value = -25281.2

# GSM!row132
$ws6.Range("H132").Value = 36396.535
$ws6.Range("I132").Value = 3392.6667
$ws6.Range("K132").Value = 10178.0001
$ws6.Range("M132").Value = -7648.000100000001

# LTW!row61
$ws7.Range("H61").Value = 2721.05
$ws7.Range("I61").Value = 1437.2858
$ws7.Range("K61").Value = 1437.2858
$ws7.Range("M61").Value = -1235.2858

# LTW!row68
$ws7.Range("H68").Value = 2741.4167
$ws7.Range("I68").Value = 2612.625
$ws7.Range("J68").Value = 2999
$ws7.Range("K68").Value = 2612.625
$ws7.Range("L68").Value = 2999
$ws7.Range("M68").Value = -1863.625
$ws7.Range("N68").Value = -4497

# LTW!row71
$ws7.Range("H71").Value = 2741.4167
$ws7.Range("I71").Value = 2612.625
$ws7.Range("J71").Value = 2999
$ws7.Range("K71").Value = 13063.125
$ws7.Range("L71").Value = 14995
$ws7.Range("M71").Value = -9319.125
$ws7.Range("N71").Value = -22483

# LTW!row93
$ws7.Range("H93").Value = 1654.2727
$ws7.Range("I93").Value = 1274.75
$ws7.Range("J93").Value = 2666.3333
$ws7.Range("K93").Value = 1274.75
$ws7.Range("L93").Value = 2666.3333
$ws7.Range("M93").Value = -26.75
$ws7.Range("N93").Value = -5162.3333

# LTW!row113
$ws7.Range("H113").Value = 2721.05
$ws7.Range("I113").Value = 1437.2858
$ws7.Range("K113").Value = 1437.2858
$ws7.Range("M113").Value = 732.7141999999999

# LTW!row122
$ws7.Range("H122").Value = 1228439.5
$ws7.Range("I122").Value = 1636402.9
$ws7.Range("K122").Value = 4909208.699999999
$ws7.Range("M122").Value = -4906758.699999999

# LTW!row132
$ws7.Range("H132").Value = 863084.2
$ws7.Range("I132").Value = 1340665.1
$ws7.Range("J132").Value = 3438.6
$ws7.Range("K132").Value = 4021995.3
$ws7.Range("L132").Value = 10315.8
$ws7.Range("M132").Value = -4019465.3
$ws7.Range("N132").Value = -15375.8

